$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Disgust" row (row 6) entirely; this shifts "Surprise" (row 7) up to row 6
$ws.Rows.Item(6).Delete()

# Update header row text
$ws.Range("B1").Value = "Hume (speech)"
$ws.Range("C1").Value = "NLP (text)"
$ws.Range("D1").Value = "Self" + [char]0x2011 + "label"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Row 2 - Anger
$ws.Range("B2").Value = 0.27
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0.08

# Row 3 - Joy
$ws.Range("B3").Value = 0.27
$ws.Range("C3").Value = 0.74
$ws.Range("D3").Value = 0.5

# Row 4 - Sadness
$ws.Range("B4").Value = 0.14
$ws.Range("C4").Value = 0.09
$ws.Range("D4").Value = 0.08

# Row 5 - Fear
$ws.Range("B5").Value = 0.18
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.08

# Row 6 - Surprise (previously row 7, shifted up after row delete)
$ws.Range("A6").Value = "Surprise"
$ws.Range("B6").Value = 0.14
$ws.Range("C6").Value = 0.17
$ws.Range("D6").Value = 0.25
